$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NATMI LR-pair metrics (Fndc5-Itgb5) with recomputed TPM-based values
# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.685737
$ws.Range("H2").Value = 2.057211
$ws.Range("I2").Value = 0.05519567570004053
$ws.Range("J2").Value = 0.05519567570004053
$ws.Range("M2").Value = 3.303267
$ws.Range("N2").Value = 9.909801000000002
$ws.Range("O2").Value = 0.03362563178859915
$ws.Range("P2").Value = 0.03362563178859915
$ws.Range("Q2").Value = 2.265172402779001
$ws.Range("R2").Value = 20.38655162501101
$ws.Range("S2").Value = 0.001855989467412493
$ws.Range("T2").Value = 0.001855989467412493

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.685737
$ws.Range("H3").Value = 2.057211
$ws.Range("I3").Value = 0.05519567570004053
$ws.Range("J3").Value = 0.05519567570004053
$ws.Range("M3").Value = 37.82684066666667
$ws.Range("O3").Value = 0.3850586149964086
$ws.Range("P3").Value = 0.3850586149964086
$ws.Range("Q3").Value = 25.939264238238
$ws.Range("R3").Value = 233.453378144142
$ws.Range("S3").Value = 0.02125357043884854
$ws.Range("T3").Value = 0.02125357043884854

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.685737
$ws.Range("H4").Value = 2.057211
$ws.Range("I4").Value = 0.05519567570004053
$ws.Range("J4").Value = 0.05519567570004053
$ws.Range("M4").Value = 9.149395999999999
$ws.Range("N4").Value = 27.448188
$ws.Range("O4").Value = 0.09313634682999644
$ws.Range("P4").Value = 0.09313634682999644
$ws.Range("Q4").Value = 6.274079364852
$ws.Range("R4").Value = 56.466714283668
$ws.Range("S4").Value = 0.005140723595514982
$ws.Range("T4").Value = 0.005140723595514982

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.685737
$ws.Range("H5").Value = 2.057211
$ws.Range("I5").Value = 0.05519567570004053
$ws.Range("J5").Value = 0.05519567570004053
$ws.Range("M5").Value = 47.95707433333333
$ws.Range("N5").Value = 143.871223
$ws.Range("O5").Value = 0.4881794063849957
$ws.Range("P5").Value = 0.4881794063849957
$ws.Range("Q5").Value = 32.885940282117
$ws.Range("R5").Value = 295.973462539053
$ws.Range("S5").Value = 0.02694539219826452
$ws.Range("T5").Value = 0.02694539219826452

# Row 6
$ws.Range("I6").Value = 0.2871009238089374
$ws.Range("J6").Value = 0.2871009238089374
$ws.Range("M6").Value = 3.303267
$ws.Range("N6").Value = 9.909801000000002
$ws.Range("O6").Value = 0.03362563178859915
$ws.Range("P6").Value = 0.03362563178859915
$ws.Range("Q6").Value = 11.782319559934
$ws.Range("R6").Value = 106.040876039406
$ws.Range("S6").Value = 0.009653949950165988
$ws.Range("T6").Value = 0.00965394995016599

# Row 7
$ws.Range("I7").Value = 0.2871009238089374
$ws.Range("J7").Value = 0.2871009238089374
$ws.Range("M7").Value = 37.82684066666667
$ws.Range("O7").Value = 0.3850586149964086
$ws.Range("P7").Value = 0.3850586149964086
$ws.Range("S7").Value = 0.1105506840860589
$ws.Range("T7").Value = 0.1105506840860589

# Row 8
$ws.Range("I8").Value = 0.2871009238089374
$ws.Range("J8").Value = 0.2871009238089374
$ws.Range("M8").Value = 9.149395999999999
$ws.Range("N8").Value = 27.448188
$ws.Range("O8").Value = 0.09313634682999644
$ws.Range("P8").Value = 0.09313634682999644
$ws.Range("Q8").Value = 32.63469391132533
$ws.Range("R8").Value = 293.712245201928
$ws.Range("S8").Value = 0.02673953121508157
$ws.Range("T8").Value = 0.02673953121508158

# Row 9
$ws.Range("I9").Value = 0.2871009238089374
$ws.Range("J9").Value = 0.2871009238089374
$ws.Range("M9").Value = 47.95707433333333
$ws.Range("N9").Value = 143.871223
$ws.Range("O9").Value = 0.4881794063849957
$ws.Range("P9").Value = 0.4881794063849957
$ws.Range("Q9").Value = 171.0565857845709
$ws.Range("R9").Value = 1539.509272061138
$ws.Range("S9").Value = 0.1401567585576309
$ws.Range("T9").Value = 0.1401567585576309

# Row 10
$ws.Range("G10").Value = 7.915626666666667
$ws.Range("H10").Value = 23.74688
$ws.Range("I10").Value = 0.6371369234209706
$ws.Range("J10").Value = 0.6371369234209707
$ws.Range("M10").Value = 3.303267
$ws.Range("N10").Value = 9.909801000000002
$ws.Range("O10").Value = 0.03362563178859915
$ws.Range("P10").Value = 0.03362563178859915
$ws.Range("Q10").Value = 26.14742835232001
$ws.Range("R10").Value = 235.32685517088
$ws.Range("S10").Value = 0.02142413158587446
$ws.Range("T10").Value = 0.02142413158587446

# Row 11
$ws.Range("G11").Value = 7.915626666666667
$ws.Range("H11").Value = 23.74688
$ws.Range("I11").Value = 0.6371369234209706
$ws.Range("J11").Value = 0.6371369234209707
$ws.Range("M11").Value = 37.82684066666667
$ws.Range("O11").Value = 0.3850586149964086
$ws.Range("P11").Value = 0.3850586149964086
$ws.Range("Q11").Value = 299.4231486968178
$ws.Range("R11").Value = 2694.80833827136
$ws.Range("S11").Value = 0.2453350612955518
$ws.Range("T11").Value = 0.2453350612955519

# Row 12
$ws.Range("G12").Value = 7.915626666666667
$ws.Range("H12").Value = 23.74688
$ws.Range("I12").Value = 0.6371369234209706
$ws.Range("J12").Value = 0.6371369234209707
$ws.Range("M12").Value = 9.149395999999999
$ws.Range("N12").Value = 27.448188
$ws.Range("O12").Value = 0.09313634682999644
$ws.Range("P12").Value = 0.09313634682999644
$ws.Range("Q12").Value = 72.42320296149333
$ws.Range("R12").Value = 651.80882665344
$ws.Range("S12").Value = 0.0593406054779324
$ws.Range("T12").Value = 0.05934060547793241

# Row 13
$ws.Range("G13").Value = 7.915626666666667
$ws.Range("H13").Value = 23.74688
$ws.Range("I13").Value = 0.6371369234209706
$ws.Range("J13").Value = 0.6371369234209707
$ws.Range("M13").Value = 47.95707433333333
$ws.Range("N13").Value = 143.871223
$ws.Range("O13").Value = 0.4881794063849957
$ws.Range("P13").Value = 0.4881794063849957
$ws.Range("Q13").Value = 379.6102964482489
$ws.Range("R13").Value = 3416.49266803424
$ws.Range("S13").Value = 0.3110371250616119
$ws.Range("T13").Value = 0.311037125061612

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.2555126666666667
$ws.Range("H14").Value = 0.7665379999999999
$ws.Range("I14").Value = 0.02056647707005147
$ws.Range("J14").Value = 0.02056647707005147
$ws.Range("M14").Value = 3.303267
$ws.Range("N14").Value = 9.909801000000002
$ws.Range("O14").Value = 0.03362563178859915
$ws.Range("P14").Value = 0.03362563178859915
$ws.Range("Q14").Value = 0.8440265598820001
$ws.Range("R14").Value = 7.596239038938001
$ws.Range("S14").Value = 0.0006915607851462185
$ws.Range("T14").Value = 0.0006915607851462185

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.2555126666666667
$ws.Range("H15").Value = 0.7665379999999999
$ws.Range("I15").Value = 0.02056647707005147
$ws.Range("J15").Value = 0.02056647707005147
$ws.Range("M15").Value = 37.82684066666667
$ws.Range("O15").Value = 0.3850586149964086
$ws.Range("P15").Value = 0.3850586149964086
$ws.Range("Q15").Value = 9.665236930315112
$ws.Range("R15").Value = 86.987132372836
$ws.Range("S15").Value = 0.007919299175949418
$ws.Range("T15").Value = 0.007919299175949418

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.2555126666666667
$ws.Range("H16").Value = 0.7665379999999999
$ws.Range("I16").Value = 0.02056647707005147
$ws.Range("J16").Value = 0.02056647707005147
$ws.Range("M16").Value = 9.149395999999999
$ws.Range("N16").Value = 27.448188
$ws.Range("O16").Value = 0.09313634682999644
$ws.Range("P16").Value = 0.09313634682999644
$ws.Range("Q16").Value = 2.337786570349333
$ws.Range("R16").Value = 21.040079133144
$ws.Range("S16").Value = 0.001915486541467483
$ws.Range("T16").Value = 0.001915486541467483

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.2555126666666667
$ws.Range("H17").Value = 0.7665379999999999
$ws.Range("I17").Value = 0.02056647707005147
$ws.Range("J17").Value = 0.02056647707005147
$ws.Range("M17").Value = 47.95707433333333
$ws.Range("N17").Value = 143.871223
$ws.Range("O17").Value = 0.4881794063849957
$ws.Range("P17").Value = 0.4881794063849957
$ws.Range("Q17").Value = 12.25363994844155
$ws.Range("R17").Value = 110.282759535974
$ws.Range("S17").Value = 0.01004013056748835
$ws.Range("T17").Value = 0.01004013056748835

Write-Host "Updated Fndc5-Itgb5 LR-pair sheet with new TPM values"
